# Chức năng Quy trình BNS + Function ExportTOExcel + MemberDetail + News
#
# Appends a new data row (row 7) to Sheet1 for a "Quy trình đào tạo"
# procedure entry. C7 ("09/09/2022") and F7 ("1233") are stored as plain
# text (shared strings) in the source workbook rather than a date serial
# or a number, so they are produced via a scratch formula cell that is
# copied and pasted back in as values/text, then cleared - this keeps the
# destination cells on the default "Normal" style (no stray NumberFormat
# bleeding into the style table), matching a plain typed-in text entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Quy trình đào tạo"

$scratch = $ws.Range("Z1")

$scratch.Formula = '="09/09/2022"'
$scratch.Copy()
$ws.Range("C7").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("D7").Value = "<p><strong>Hello</strong></p>"
$ws.Range("E7").Value = "Chưa duyệt"

$scratch.Formula = '="1233"'
$scratch.Copy()
$ws.Range("F7").PasteSpecial(-4163)
$scratch.Clear()
